$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ A = "20160406_104427"; B = 1112.233; C = 'convert to lower, trim "space" and ",", remove multiple spaces, convert unicode to ascii'; D = '7 features: #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_2, first_character_type, #"space"'; E = "Neuron Network"; F = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 1000"; G = 0.998666666666667; H = 0.871287128712871; I = "0 filters: "; J = 0 },
    @{ A = "20160406_110259"; B = 1146.936; C = 'convert to lower, trim "space" and ",", remove multiple spaces, convert unicode to ascii'; D = '7 features: #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_2, first_character_type, #"space"'; E = "Neuron Network"; F = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 1000"; G = 0.999333333333333; H = 0.867986798679868; I = "0 filters: "; J = 0.0169491525423729 },
    @{ A = "20160406_112206"; B = 1186.473; C = 'convert to lower, trim "space" and ",", remove multiple spaces, convert unicode to ascii'; D = '7 features: #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_2, first_character_type, #"space"'; E = "Neuron Network"; F = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 1000"; G = 1; H = 0.877887788778878; I = "0 filters: "; J = 0.0161290322580645 },
    @{ A = "20160406_114152"; B = 1200.741; C = 'convert to lower, trim "space" and ",", remove multiple spaces, convert unicode to ascii'; D = '7 features: #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_2, first_character_type, #"space"'; E = "Neuron Network"; F = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 1000"; G = 1; H = 0.867986798679868; I = "0 filters: "; J = 0.0338983050847458 },
    @{ A = "20160406_120153"; B = 1262.807; C = 'convert to lower, trim "space" and ",", remove multiple spaces, convert unicode to ascii'; D = '7 features: #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_2, first_character_type, #"space"'; E = "Neuron Network"; F = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 1000"; G = 0.998666666666667; H = 0.884488448844885; I = "0 filters: "; J = 0.015625 },
    @{ A = "20160406_133609"; B = 2579.909; C = 'trim "space" and ",", remove multiple spaces, convert unicode to ascii, convert to lower'; D = '7 features: #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_2, first_character_type, #"space"'; E = "Neuron Network"; F = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 2000"; G = 1; H = 0.867986798679868; I = "0 filters: "; J = 0 },
    @{ A = "20160406_141909"; B = 2681.746; C = 'trim "space" and ",", remove multiple spaces, convert unicode to ascii, convert to lower'; D = '7 features: #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_2, first_character_type, #"space"'; E = "Neuron Network"; F = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 2000"; G = 0.999333333333333; H = 0.864686468646865; I = "0 filters: "; J = 0 },
    @{ A = "20160406_150350"; B = 1682.843; C = 'trim "space" and ",", remove multiple spaces, convert unicode to ascii, convert to lower'; D = '7 features: #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_2, first_character_type, #"space"'; E = "Neuron Network"; F = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 2000"; G = 0.999333333333333; H = 0.867986798679868; I = "0 filters: "; J = 0 },
    @{ A = "20160406_153153"; B = 1605.124; C = 'trim "space" and ",", remove multiple spaces, convert unicode to ascii, convert to lower'; D = '7 features: #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_2, first_character_type, #"space"'; E = "Neuron Network"; F = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 2000"; G = 0.999333333333333; H = 0.867986798679868; I = "0 filters: "; J = 0 },
    @{ A = "20160406_155838"; B = 1514.518; C = 'trim "space" and ",", remove multiple spaces, convert unicode to ascii, convert to lower'; D = '7 features: #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_2, first_character_type, #"space"'; E = "Neuron Network"; F = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 2000"; G = 1; H = 0.867986798679868; I = "0 filters: "; J = 0 }
)

$startRow = 12
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
}
